# Add monitor site IWIP
# - Fill in the NE (L2SW) model column with more specific hardware names
# - Replace the placeholder NMS credentials with the shared username_nms /
#   password_nms account across the L2SW and METRO device blocks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("NE") - L2SW device model detail
$ws.Range("C2:C4").Value = "L2SW FH S5800v3"
$ws.Range("C5:C27").Value = "L2SW FH S5800v2"

# Columns F/G ("Username NE" / "Password NE") - shared NMS credentials
$ws.Range("F5:F35").Value = "username_nms"
$ws.Range("G5:G35").Value = "password_nms"

# Reflect where the user's selection ended up after making the edit
$ws.Range("F9").Select()

$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null
